# Updated cryptos list on Sun Apr 23 22:36:01 UTC 2023 with GitHub Actions
#
# This script refreshes the Price (column D) and Volume(1h) (column E)
# values for the crypto-ranking sheet, and fixes the Frax / FraxShare
# rows (40-41) which had been swapped so that the coin name / link now
# line up with the correct price & volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    # Column D prices look like "1.017" / "27.916.05" - plain ".Value"
    # assignment lets Excel's COM layer auto-coerce single-dot numeric
    # looking strings into real numbers, which would corrupt the text.
    # Forcing the cell to Text format first keeps it a string, matching
    # the original inline-string content.
    $ws.Range($CellRef).NumberFormat = "@"
    $ws.Range($CellRef).Value = $Text
}

# D2/D3 (and a few others below) already fail Excel's numeric-literal
# check because of the repeated "." separators, so a plain assignment
# is kept as text without needing the NumberFormat coercion.
$ws.Range("D2").Value = "27.949.68"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.887.69"
$ws.Range("E3").Value = "  +0.41%  "

Set-PriceText "D4" "1.017"
$ws.Range("E4").Value = "  +1.39%  "

Set-PriceText "D5" "335.83"
$ws.Range("E5").Value = "  +0.69%  "

Set-PriceText "D6" "1.016"
$ws.Range("E6").Value = "  +1.34%  "

Set-PriceText "D7" "0.4693"
$ws.Range("E7").Value = "  -0.94%  "

Set-PriceText "D8" "0.3913"
$ws.Range("E8").Value = "  -1.64%  "

Set-PriceText "D9" "46.95"
$ws.Range("E9").Value = "  -3.31%  "

Set-PriceText "D10" "0.07970"
$ws.Range("E10").Value = "  -0.92%  "

Set-PriceText "D11" "1.014"
$ws.Range("E11").Value = "  -1.21%  "

Set-PriceText "D12" "21.77"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "1.885.72"
$ws.Range("E13").Value = "  -1.74%  "

Set-PriceText "D14" "5.964"
$ws.Range("E14").Value = "  -0.11%  "

Set-PriceText "D15" "7.132"
$ws.Range("E15").Value = "  -0.81%  "

Set-PriceText "D16" "1.018"
$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("E17").Value = "  +2.34%  "

Set-PriceText "D18" "87.55"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("E19").Value = "  -0.51%  "

Set-PriceText "D20" "17.05"
$ws.Range("E20").Value = "  -1.75%  "

Set-PriceText "D21" "1.016"
$ws.Range("E21").Value = "  +1.30%  "

$ws.Range("D22").Value = "27.944.11"
$ws.Range("E22").Value = "  -0.12%  "

Set-PriceText "D23" "5.483"
$ws.Range("E23").Value = "  -0.50%  "

Set-PriceText "D24" "10.95"
$ws.Range("E24").Value = "  -0.96%  "

Set-PriceText "D25" "2.362"
$ws.Range("E25").Value = "  +2.60%  "

$ws.Range("D26").Value = "2.104.40"
$ws.Range("E26").Value = "  -1.60%  "

Set-PriceText "D27" "159.72"
$ws.Range("E27").Value = "  +1.54%  "

Set-PriceText "D28" "20.02"
$ws.Range("E28").Value = "  -1.24%  "

Set-PriceText "D29" "2.092"
$ws.Range("E29").Value = "  -0.75%  "

Set-PriceText "D30" "5.491"
$ws.Range("E30").Value = "  -2.10%  "

Set-PriceText "D31" "121.24"
$ws.Range("E31").Value = "  -1.31%  "

Set-PriceText "D32" "0.09564"
$ws.Range("E32").Value = "  -0.10%  "

Set-PriceText "D33" "0.9601"
$ws.Range("E33").Value = "  -1.92%  "

Set-PriceText "D34" "3.653"
$ws.Range("E34").Value = "  +0.41%  "

Set-PriceText "D35" "5.342"
$ws.Range("E35").Value = "  +0.38%  "

Set-PriceText "D36" "1.354"
$ws.Range("E36").Value = "  -7.74%  "

Set-PriceText "D37" "0.06124"
$ws.Range("E37").Value = "  +0.07%  "

Set-PriceText "D38" "0.02248"
$ws.Range("E38").Value = "  -0.67%  "

Set-PriceText "D39" "1.209"
$ws.Range("E39").Value = "  -2.07%  "

# Rows 40 and 41 were swapped: row 40 now holds FraxShare (FXS) data and
# row 41 now holds Frax (FRAX) data.
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-PriceText "D40" "8.162"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-PriceText "D41" "1.016"
$ws.Range("E41").Value = "  +1.37%  "

Set-PriceText "D42" "0.5935"
$ws.Range("E42").Value = "  -1.68%  "

Set-PriceText "D43" "0.1897"
$ws.Range("E43").Value = "  -0.80%  "

Set-PriceText "D44" "10.27"
$ws.Range("E44").Value = "  -0.58%  "

Set-PriceText "D45" "1.269"
$ws.Range("E45").Value = "  +1.78%  "

Set-PriceText "D46" "0.5674"
$ws.Range("E46").Value = "  -0.81%  "

Set-PriceText "D47" "12.18"
$ws.Range("E47").Value = "  -0.75%  "

Set-PriceText "D48" "3.406"
$ws.Range("E48").Value = "  -0.36%  "

Set-PriceText "D49" "1.936"
$ws.Range("E49").Value = "  -0.40%  "

Set-PriceText "D50" "0.06853"
$ws.Range("E50").Value = "  +0.33%  "

Set-PriceText "D51" "114.01"
$ws.Range("E51").Value = "  +0.10%  "
